{"js": "// Regional settings fix (Spain): wrap the numeric merge-fields used in the\n// sales-order table (unit price, line amount and grand total) with the\n// docxtemplater `formatN(2)` filter so that numbers render using the\n// Spanish locale's numeric formatting (2 decimal places).\n//\n//   { d.order.salesOrderDetails[i].unitPrice }   -> { ... .unitPrice:formatN(2) }\n//   { d.order.salesOrderDetails[i].amount }      -> { ... .amount:formatN(2) }\n//   { d.total }                                  -> { d.total:formatN(2) }\n\nconst body = context.document.body;\n\nasync function appendFormatN(searchText) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(\":formatN(2)\", Word.InsertLocation.after);\n  }\n  await context.sync();\n}\n\n// Unit price and line amount appear twice each (loop body rendered for the\n// two visible rows `i` and `i+1` of the sales-order template).\nawait appendFormatN(\"unitPrice\");\nawait appendFormatN(\"amount\");\n\n// Grand total (bold, single occurrence).\nawait appendFormatN(\"d.total\");\n", "ps1": "# Regional settings fix (Spain): wrap the numeric merge-fields used in the\n# sales-order table (unit price, line amount and grand total) with the\n# docxtemplater `formatN(2)` filter so that numbers render using the\n# Spanish locale's numeric formatting (2 decimal places).\n#\n#   { d.order.salesOrderDetails[i].unitPrice }   -> { ... .unitPrice:formatN(2) }\n#   { d.order.salesOrderDetails[i].amount }      -> { ... .amount:formatN(2) }\n#   { d.total }                                  -> { d.total:formatN(2) }\n\n$d = $word.ActiveDocument\n\nfunction Append-FormatN($searchText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # wdFindContinue = 1, wdReplaceAll = 2 ; MatchWholeWord = $true, MatchCase = $true\n    $find.Execute($searchText, $true, $true, $false, $false, $false, $true, 1, $false, ($searchText + \":formatN(2)\"), 2) | Out-Null\n}\n\n# Unit price and line amount appear twice each (loop body rendered for the\n# two visible rows `i` and `i+1` of the sales-order template).\nAppend-FormatN \"unitPrice\"\nAppend-FormatN \"amount\"\n\n# Grand total (bold, single occurrence).\nAppend-FormatN \"d.total\"\n"}
